$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-7) hold weekly price records for "Caqui" that got re-sorted.
# Only columns D (Fecha), K (Variedad), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen) and S (Precio $/Kg) differ between rows, so we only need to
# rewrite those columns for each row according to the new order.

$rows = @{
    2 = @{ D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    3 = @{ D = 44355; K = "Mankaki"; L = "Segunda"; M = 270; N = 20000; O = 21000; P = 20500; R = "Región Metropolitana"; S = 1139 }
    4 = @{ D = 44699; K = "Mankaki"; L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins"; S = 1639 }
    5 = @{ D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 }
    6 = @{ D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 }
    7 = @{ D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
